# Processing Uploaded file and Rspec coverage
#
# The workbook originally has 6 sheets: Sheet1 (empty placeholder) and
# Sheet2..Sheet6 (each containing a survey question's grouped data).
# This script removes the empty leading sheet and renames the remaining
# five data sheets to "Question 11".."Question 15".

$wb = $excel.ActiveWorkbook

# Delete the empty first worksheet (old "Sheet1").
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet1").Delete()
$excel.DisplayAlerts = $true

# Rename the remaining sheets (old Sheet2..Sheet6) to Question 11..Question 15.
$wb.Worksheets.Item("Sheet2").Name = "Question 11"
$wb.Worksheets.Item("Sheet3").Name = "Question 12"
$wb.Worksheets.Item("Sheet4").Name = "Question 13"
$wb.Worksheets.Item("Sheet5").Name = "Question 14"
$wb.Worksheets.Item("Sheet6").Name = "Question 15"
